$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = '28.369.29'
$ws.Range("E2").Value2 = '  -0.45%  '

$ws.Range("D3").Value2 = '1.564.28'
$ws.Range("E3").Value2 = '  -0.19%  '

$ws.Range("E4").Value2 = '  +0.00%  '

$ws.Range("D5").Value2 = '211.07'
$ws.Range("E5").Value2 = '  -0.33%  '

$ws.Range("E6").Value2 = '  -0.64%  '

$ws.Range("E7").Value2 = '  +0.01%  '

$ws.Range("D8").Value2 = '44.54'
$ws.Range("E8").Value2 = '  -3.66%  '

$ws.Range("D9").Value2 = '23.56'
$ws.Range("E9").Value2 = '  -2.03%  '

$ws.Range("E10").Value2 = '  -1.42%  '

$ws.Range("D11").Value2 = '0.0588'
$ws.Range("E11").Value2 = '  -0.67%  '

$ws.Range("E12").Value2 = '  +0.78%  '

$ws.Range("D13").Value2 = '1.787.91'
$ws.Range("E13").Value2 = '  -0.16%  '

$ws.Range("D14").Value2 = '1.569.45'
$ws.Range("E14").Value2 = '  +0.17%  '

$ws.Range("E15").Value2 = '  -0.39%  '

$ws.Range("D16").Value2 = '28.360.78'

$ws.Range("D17").Value2 = '0.512'
$ws.Range("E17").Value2 = '  -1.46%  '

$ws.Range("D18").Value2 = '60.46'
$ws.Range("E18").Value2 = '  -2.85%  '

$ws.Range("D19").Value2 = '228.16'
$ws.Range("E19").Value2 = '  +0.06%  '

$ws.Range("E20").Value2 = '  +0.20%  '

$ws.Range("E21").Value2 = '  -1.87%  '

$ws.Range("E22").Value2 = '  +0.00%  '

$ws.Range("E23").Value2 = '  +1.27%  '

$ws.Range("E24").Value2 = '  -2.06%  '

$ws.Range("E25").Value2 = '  -1.29%  '

$ws.Range("D26").Value2 = '150.28'
$ws.Range("E26").Value2 = '  -0.30%  '

$ws.Range("D27").Value2 = '14.88'
$ws.Range("E27").Value2 = '  -0.82%  '

$ws.Range("E28").Value2 = '  +0.20%  '

$ws.Range("E29").Value2 = '  -2.12%  '

$ws.Range("E30").Value2 = '  +0.00%  '

$ws.Range("E31").Value2 = '  +1.80%  '

$ws.Range("E32").Value2 = '  -4.04%  '

$ws.Range("E33").Value2 = '  -1.20%  '

$ws.Range("E34").Value2 = '  +0.00%  '

$ws.Range("D35").Value2 = '1.385.20'
$ws.Range("E35").Value2 = '  -0.60%  '

$ws.Range("E36").Value2 = '  +1.79%  '

$ws.Range("E37").Value2 = '  -3.36%  '

$ws.Range("E38").Value2 = '  -0.19%  '

$ws.Range("D39").Value2 = '2.64'
$ws.Range("E39").Value2 = '  +2.94%  '

$ws.Range("E40").Value2 = '  -2.00%  '

$ws.Range("D41").Value2 = '0.519'
$ws.Range("E41").Value2 = '  -3.23%  '

$ws.Range("E42").Value2 = '  +2.87%  '

$ws.Range("E43").Value2 = '  -0.02%  '

$ws.Range("E44").Value2 = '  -0.27%  '

$ws.Range("D45").Value2 = '0.0468'
$ws.Range("E45").Value2 = '  -2.26%  '

$ws.Range("E46").Value2 = '  -2.89%  '

$ws.Range("B47").Value2 = 'WEMIXToken'
$ws.Range("C47").Value2 = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D47").Value2 = '0.922'
$ws.Range("E47").Value2 = '  -5.40%  '

$ws.Range("B48").Value2 = 'Aave'
$ws.Range("C48").Value2 = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D48").Value2 = '62.16'
$ws.Range("E48").Value2 = '  -1.15%  '

$ws.Range("D49").Value2 = '1.700.74'
$ws.Range("E49").Value2 = '  -0.14%  '

$ws.Range("D50").Value2 = '85.38'
$ws.Range("E50").Value2 = '  -0.64%  '

$ws.Range("E51").Value2 = '  -2.19%  '
